# Updated cryptos list on Mon Feb 12 09:57:00 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'47.826.99"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "'2.474.82"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'315.55"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("D6").Value = "'104.66"
$ws.Range("E6").Value = "  -4.86%  "
$ws.Range("D7").Value = "'0.517"
$ws.Range("E7").Value = "  -3.17%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  -4.10%  "
$ws.Range("D10").Value = "'38.65"
$ws.Range("E10").Value = "  -5.18%  "
$ws.Range("D11").Value = "'19.99"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("D12").Value = "'0.0797"
$ws.Range("E12").Value = "  -3.58%  "
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "'7.03"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("D15").Value = "'2.870.57"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "'2.485.45"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "'0.822"
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("D18").Value = "'47.790.32"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "'2.95"
$ws.Range("E19").Value = "  +8.53%  "
$ws.Range("D20").Value = "'12.60"
$ws.Range("E20").Value = "  -4.77%  "
$ws.Range("D21").Value = "'6.50"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").Value = "'0.0₃0925"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").Value = "'272.34"
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("D24").Value = "'70.49"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").Value = "'2.49"
$ws.Range("E25").Value = "  -3.70%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'25.48"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("D28").Value = "'2.19"
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("D29").Value = "'9.61"
$ws.Range("E29").Value = "  -5.41%  "
$ws.Range("D30").Value = "'0.137"
$ws.Range("E30").Value = "  -5.34%  "
$ws.Range("D31").Value = "'34.21"
$ws.Range("E31").Value = "  -6.06%  "
$ws.Range("D32").Value = "'49.22"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "'18.84"
$ws.Range("E34").Value = "  -5.52%  "
$ws.Range("D35").Value = "'5.22"
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("D36").Value = "'0.0766"
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("D38").Value = "'4.49"
$ws.Range("E38").Value = "  -4.86%  "
$ws.Range("D39").Value = "'2.84"
$ws.Range("E39").Value = "  -5.37%  "
$ws.Range("D40").Value = "'122.24"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").Value = "'0.110"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "'0.0300"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "'1.990.54"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").Value = "'3.13"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "'1.89"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").Value = "'8.86"
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("D50").Value = "'5.12"
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("D51").Value = "'77.96"
$ws.Range("E51").Value = "  -1.66%  "
